$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.70"
$ws.Range("E2").Value = "'-1.71%"
$ws.Range("D3").Value = "'37.32"
$ws.Range("E3").Value = "'7.14%"
$ws.Range("D4").Value = "'4.994"
$ws.Range("E4").Value = "'-3.44%"
$ws.Range("D5").Value = "'0.07805"
$ws.Range("E5").Value = "'0.45%"
$ws.Range("D6").Value = "'2.204"
$ws.Range("E6").Value = "'-6.75%"
$ws.Range("D7").Value = "'8.018"
$ws.Range("E7").Value = "'-0.15%"
$ws.Range("D8").Value = "'4.029"
$ws.Range("E8").Value = "'1.66%"
$ws.Range("D9").Value = "'0.9147"
$ws.Range("E9").Value = "'-1.47%"
$ws.Range("D10").Value = "'0.09690"
$ws.Range("E10").Value = "'-4.81%"
$ws.Range("D11").Value = "'0.1889"
$ws.Range("E11").Value = "'2.82%"
$ws.Range("D12").Value = "'0.08721"
$ws.Range("E12").Value = "'0.16%"
$ws.Range("E13").Value = "'3.72%"
$ws.Range("D14").Value = "'0.09948"
$ws.Range("E14").Value = "'0.68%"
$ws.Range("D15").Value = "'0.001488"
$ws.Range("E15").Value = "'0.27%"
$ws.Range("D16").Value = "'0.005634"
$ws.Range("E16").Value = "'-0.54%"
$ws.Range("D17").Value = "'3.458"
$ws.Range("E17").Value = "'-1.27%"
$ws.Range("E18").Value = "'12.72%"
$ws.Range("E19").Value = "'1.89%"
$ws.Range("D20").Value = "'0.1276"
$ws.Range("D21").Value = "'4.777"
$ws.Range("E21").Value = "'5.59%"
$ws.Range("E22").Value = "'0.24%"
$ws.Range("D23").Value = "'0.04635"
$ws.Range("E23").Value = "'1.03%"
$ws.Range("D24").Value = "'0.001232"
$ws.Range("E24").Value = "'1.41%"
$ws.Range("D25").Value = "'0.004790"
$ws.Range("E25").Value = "'6.93%"
$ws.Range("E26").Value = "'-7.08%"
$ws.Range("E27").Value = "'40.23%"
$ws.Range("D39").Value = "'0.01754"
$ws.Range("E39").Value = "'-2.31%"
$ws.Range("D40").Value = "'0.04741"
$ws.Range("E40").Value = "'-1.37%"
$ws.Range("D41").Value = "'0.008076"
$ws.Range("E41").Value = "'4.65%"
$ws.Range("D42").Value = "'0.1389"
$ws.Range("E42").Value = "'-2.13%"
$ws.Range("D43").Value = "'0.007690"
$ws.Range("D44").Value = "'0.002234"
$ws.Range("E44").Value = "'1.30%"
$ws.Range("D45").Value = "'0.01041"
$ws.Range("E45").Value = "'8.52%"
$ws.Range("D46").Value = "'0.00006015"
$ws.Range("E46").Value = "'0.59%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.50%"
$ws.Range("D48").Value = "'7.828"
$ws.Range("E48").Value = "'186.48%"
$ws.Range("E49").Value = "'0.41%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.50%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.50%"
